$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 281.25
$ws.Range("I33").Value = 281.25
$ws.Range("K33").Value = 281.25
$ws.Range("M33").Value = -52.25

$ws.Range("H54").Value = 9538
$ws.Range("I54").Value = 9538
$ws.Range("K54").Value = 9538
$ws.Range("M54").Value = -9052

$ws.Range("H86").Value = 2611.111
$ws.Range("J86").Value = 2875
$ws.Range("L86").Value = 2875
$ws.Range("N86").Value = -5121

$ws.Range("H89").Value = 2611.111
$ws.Range("J89").Value = 2875
$ws.Range("L89").Value = 14375
$ws.Range("N89").Value = -25607

$ws.Range("H131").Value = 6317
$ws.Range("I131").Value = 1162.7142
$ws.Range("J131").Value = 9203.4
$ws.Range("K131").Value = 3488.1426
$ws.Range("L131").Value = 27610.2
$ws.Range("M131").Value = 1551.8574
$ws.Range("N131").Value = -37690.2

$ws.Range("H136").Value = 76000
$ws.Range("J136").Value = 76000
$ws.Range("L136").Value = 76000
$ws.Range("N136").Value = -86200

$ws.Range("H141").Value = 7310.048
$ws.Range("I141").Value = 2423.5881
$ws.Range("J141").Value = 28077.5
$ws.Range("K141").Value = 7270.7643
$ws.Range("L141").Value = 84232.5
$ws.Range("M141").Value = -2090.7643
$ws.Range("N141").Value = -94592.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10941.529
$ws.Range("I32").Value = 13391.519
$ws.Range("K32").Value = 13391.519
$ws.Range("M32").Value = -13104.519

$ws.Range("H88").Value = 2897.4443
$ws.Range("J88").Value = 3039.5715
$ws.Range("L88").Value = 3039.5715
$ws.Range("N88").Value = -3851.5715

$ws.Range("H91").Value = 2897.4443
$ws.Range("J91").Value = 3039.5715
$ws.Range("L91").Value = 3039.5715
$ws.Range("N91").Value = -5847.5715

$ws.Range("H122").Value = 2987.6
$ws.Range("I122").Value = 2987.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8962.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6512.799999999999
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229

$ws.Range("H124").Value = 29714.5
$ws.Range("J124").Value = 29714.5
$ws.Range("L124").Value = 29714.5
$ws.Range("N124").Value = -39534.5

$ws.Range("H125").Value = 48476.668
$ws.Range("J125").Value = 48476.668
$ws.Range("L125").Value = 48476.668
$ws.Range("N125").Value = -58316.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 25442.523
$ws.Range("J52").Value = 25442.523
$ws.Range("L52").Value = 25442.523
$ws.Range("N52").Value = -25968.523

$ws.Range("H121").Value = 25442.523
$ws.Range("J121").Value = 25442.523
$ws.Range("L121").Value = 25442.523
$ws.Range("N121").Value = -28936.523

$ws.Range("H124").Value = 40780
$ws.Range("J124").Value = 40780
$ws.Range("L124").Value = 40780
$ws.Range("N124").Value = -50600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2796.9487
$ws.Range("I31").Value = 1507.0741
$ws.Range("J31").Value = 5699.1665
$ws.Range("K31").Value = 1507.0741
$ws.Range("L31").Value = 5699.1665
$ws.Range("M31").Value = -1212.0741
$ws.Range("N31").Value = -6289.1665

$ws.Range("H34").Value = 2796.9487
$ws.Range("I34").Value = 1507.0741
$ws.Range("J34").Value = 5699.1665
$ws.Range("K34").Value = 1507.0741
$ws.Range("L34").Value = 5699.1665
$ws.Range("M34").Value = -1305.0741
$ws.Range("N34").Value = -6103.1665

$ws.Range("H58").Value = 2096.5557
$ws.Range("I58").Value = 1954.8235
$ws.Range("J58").Value = 2337.5
$ws.Range("K58").Value = 1954.8235
$ws.Range("L58").Value = 2337.5
$ws.Range("M58").Value = -1751.8235
$ws.Range("N58").Value = -2743.5

$ws.Range("H99").Value = 2181.4
$ws.Range("I99").Value = 2628.5715
$ws.Range("J99").Value = 1138
$ws.Range("K99").Value = 2628.5715
$ws.Range("L99").Value = 1138
$ws.Range("M99").Value = -1130.5715
$ws.Range("N99").Value = -4134

$ws.Range("H126").Value = 2181.4
$ws.Range("I126").Value = 2628.5715
$ws.Range("J126").Value = 1138
$ws.Range("K126").Value = 7885.7145
$ws.Range("L126").Value = 3414
$ws.Range("M126").Value = -5415.7145
$ws.Range("N126").Value = -8354

$ws.Range("H136").Value = 2096.5557
$ws.Range("I136").Value = 1954.8235
$ws.Range("J136").Value = 2337.5
$ws.Range("K136").Value = 5864.470499999999
$ws.Range("L136").Value = 7012.5
$ws.Range("M136").Value = -3314.470499999999
$ws.Range("N136").Value = -12112.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 5500
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 6000
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = -4748
$ws.Range("N33").Value = -6504

$ws.Range("H102").Value = 2116.4092
$ws.Range("I102").Value = 1979.6154
$ws.Range("J102").Value = 2314
$ws.Range("K102").Value = 1979.6154
$ws.Range("L102").Value = 2314
$ws.Range("M102").Value = -357.6153999999999
$ws.Range("N102").Value = -5558

$ws.Range("H109").Value = 17785
$ws.Range("J109").Value = 17785
$ws.Range("L109").Value = 17785
$ws.Range("N109").Value = -19865

$ws.Range("H122").Value = 2614.5
$ws.Range("I122").Value = 2038.6875
$ws.Range("J122").Value = 3535.8
$ws.Range("K122").Value = 6116.0625
$ws.Range("L122").Value = 10607.4
$ws.Range("M122").Value = -3666.0625
$ws.Range("N122").Value = -15507.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 11000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H127").Value = 79333.336
$ws.Range("J127").Value = 79333.336
$ws.Range("L127").Value = 79333.336
$ws.Range("N127").Value = -89253.336

$ws.Range("H132").Value = 5399.316
$ws.Range("I132").Value = 6460
$ws.Range("J132").Value = 4220.778
$ws.Range("K132").Value = 19380
$ws.Range("L132").Value = 12662.334
$ws.Range("M132").Value = -16850
$ws.Range("N132").Value = -17722.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 39533.332
$ws.Range("J40").Value = 39533.332
$ws.Range("L40").Value = 39533.332
$ws.Range("N40").Value = -39831.332

$ws.Range("H123").Value = 24104.314
$ws.Range("J123").Value = 24104.314
$ws.Range("L123").Value = 24104.314
$ws.Range("N123").Value = -33904.314

$ws.Range("H126").Value = 4263.778
$ws.Range("I126").Value = 4554.3335
$ws.Range("J126").Value = 3682.6667
$ws.Range("K126").Value = 13663.0005
$ws.Range("L126").Value = 11048.0001
$ws.Range("M126").Value = -11193.0005
$ws.Range("N126").Value = -15988.0001
